# Added Log Filter Buttons
#
# 1) Insert a new bold "J1.0 Order Receive: " paragraph right after the
#    existing "J0.5 Relay:" paragraph (and before the blank paragraph
#    that used to directly follow it).
# 2) Because the new paragraph pushes everything else down by one
#    paragraph, Word's cached pagination marker (<w:lastRenderedPageBreak/>)
#    shifts from the "J13.3 Move Mission Start:" run to the
#    "J13.2 Follow Mission Start:" run. Reproduce that move explicitly.

$d = $word.ActiveDocument

$wdCollapseEnd = 0

$pkgOpen = "<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'><pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'><pkg:xmlData><w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:body>"
$pkgClose = "</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"

# ---------------------------------------------------------------------
# Step 1: locate the "J0.5 Relay:" paragraph and insert the new
# "J1.0 Order Receive: " paragraph right after it.
# ---------------------------------------------------------------------
$relayPara = $null
foreach ($p in $d.Paragraphs) {
    $txt = $p.Range.Text
    $txt = $txt.TrimEnd([char]13, [char]7)
    if ($txt -eq "J0.5 Relay:") {
        $relayPara = $p
        break
    }
}

$insertRange = $relayPara.Range.Duplicate
$insertRange.Collapse($wdCollapseEnd)
$insertRange.InsertParagraphAfter()

$newHeadingPara = $relayPara.Next()
$newHeadingRange = $newHeadingPara.Range
$newHeadingRange.Text = "J1.0 Order Receive: "
$newHeadingRange.Font.Name = "Arial"
$newHeadingRange.Font.Bold = $true
$newHeadingRange.Font.Size = 12

# ---------------------------------------------------------------------
# Step 2: move <w:lastRenderedPageBreak/> from the "J13.3 Move Mission
# Start:" run to the "J13.2 Follow Mission Start:" run.
# ---------------------------------------------------------------------

$runRPrXml = "<w:rPr><w:rFonts w:ascii='Arial' w:hAnsi='Arial' w:cs='Arial'/><w:b/><w:bCs/><w:sz w:val='24'/><w:szCs w:val='24'/><w:lang w:val='en-US'/></w:rPr>"

function Replace-FirstCharWithRun($doc, $paragraph, $innerXml) {
    $full = $paragraph.Range.Duplicate
    $charRange = $doc.Range($full.Start, $full.Start + 1)
    $charRange.Delete()
    $insPoint = $doc.Range($full.Start, $full.Start)
    $xml = $pkgOpen + "<w:p>" + $innerXml + "</w:p>" + $pkgClose
    $insPoint.InsertXML($xml)
}

$para132 = $null
$para133 = $null
foreach ($p in $d.Paragraphs) {
    $txt = $p.Range.Text
    if ($txt -like "*J13.2 Follow Mission Start:*") {
        $para132 = $p
    }
    if ($txt -like "*J13.3 Move Mission Start:*") {
        $para133 = $p
    }
}

# Remove the marker from the "J13.3" run (rebuild the leading "J" run
# without <w:lastRenderedPageBreak/>).
$innerNoBreak = "<w:r w:rsidRPr='00B03554'>" + $runRPrXml + "<w:t>J</w:t></w:r>"
Replace-FirstCharWithRun $d $para133 $innerNoBreak

# Add the marker to the "J13.2" run (rebuild the leading "J" run with
# <w:lastRenderedPageBreak/> inserted before the text).
$innerWithBreak = "<w:r w:rsidRPr='00B03554'>" + $runRPrXml + "<w:lastRenderedPageBreak/><w:t>J</w:t></w:r>"
Replace-FirstCharWithRun $d $para132 $innerWithBreak
